# Add a new "References & Links" slide (Title and Content layout) at the
# end of the deck, containing the reference book titles plus a hyperlinked
# GitHub URL - matching the "Add placeholder for links." commit.

$p = $ppt.ActivePresentation

# Layout 2 == "Title and Content" (same layout used by the other
# text-only content slides in this deck, e.g. slide 9).
$s = $p.Slides.Add($p.Slides.Count + 1, 2)

$url = "https://github.com/crsdrw/patterns/tree/factory_method"

# --- Title placeholder ---------------------------------------------------
$title = $s.Shapes.Item(1).TextFrame.TextRange
# Set the language before typing so it is stamped on every run created by
# the subsequent .Text assignment.
$title.LanguageID = "en-GB"
$title.Text = "References & Links"

# --- Body / content placeholder ------------------------------------------
$body = $s.Shapes.Item(2).TextFrame.TextRange
$body.LanguageID = "en-GB"
$body.Text = "Head First Design Patterns" + "`r" + "Design Patterns" + "`r" + $url

# Turn the 3rd paragraph (the bare URL) into a hyperlink. PowerPoint's
# AutoFormat splits an auto-detected URL into separate runs around the
# "://" - reproduce that by hyperlinking three character sub-ranges.
$linkPara = $body.Paragraphs(3)
$schemeLen = 5                               # "https"
$sepLen = 3                                  # "://"
$restLen = $linkPara.Length - $schemeLen - $sepLen

$linkPara.Characters(1, $schemeLen).ActionSettings(1).Hyperlink.Address = $url
$linkPara.Characters($schemeLen + 1, $sepLen).ActionSettings(1).Hyperlink.Address = $url
$linkPara.Characters($schemeLen + $sepLen + 1, $restLen).ActionSettings(1).Hyperlink.Address = $url
